$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume table (scraped data updated by the bot).
# All D/E/B/C cells in this sheet are stored as plain text (inline strings), including
# numeric-looking prices such as '1.00' or '305.48'. For values that Excel would
# otherwise auto-convert to a true number, we force text entry with a leading
# apostrophe and then strip the resulting quote-prefix style so formatting stays default.

# Row 2
$ws.Range("D2").Value2 = '42.847.56'
$ws.Range("E2").Value2 = '  +0.10%  '

# Row 3
$ws.Range("D3").Value2 = '2.533.74'
$ws.Range("E3").Value2 = '  -1.38%  '

# Row 4
$ws.Range("D4").Value2 = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = '  -0.07%  '

# Row 5
$ws.Range("D5").Value2 = "'305.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +1.04%  '

# Row 6
$ws.Range("D6").Value2 = "'98.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  +6.33%  '

# Row 7
$ws.Range("D7").Value2 = "'0.583"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = '  +1.69%  '

# Row 8
$ws.Range("E8").Value2 = '  +0.05%  '

# Row 9
$ws.Range("D9").Value2 = "'0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = '  +0.23%  '

# Row 10
$ws.Range("E10").Value2 = '  +2.63%  '

# Row 11
$ws.Range("D11").Value2 = "'0.0814"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = '  +0.53%  '

# Row 12
$ws.Range("D12").Value2 = "'7.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = '  +0.91%  '

# Row 13
$ws.Range("E13").Value2 = '  -0.84%  '

# Row 14
$ws.Range("D14").Value2 = '2.921.81'
$ws.Range("E14").Value2 = '  -1.45%  '

# Row 15
$ws.Range("D15").Value2 = '2.543.55'
$ws.Range("E15").Value2 = '  -1.99%  '

# Row 16
$ws.Range("D16").Value2 = "'15.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = '  +6.23%  '

# Row 17
$ws.Range("D17").Value2 = "'0.868"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = '  -1.39%  '

# Row 18
$ws.Range("D18").Value2 = '42.860.85'
$ws.Range("E18").Value2 = '  -0.03%  '

# Row 19
$ws.Range("D19").Value2 = "'12.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  +1.87%  '

# Row 20
$ws.Range("D20").Value2 = '0.0₃0982'
$ws.Range("E20").Value2 = '  -0.72%  '

# Row 21
$ws.Range("D21").Value2 = "'6.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  -2.09%  '

# Row 22
$ws.Range("D22").Value2 = "'71.50"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").Value2 = "'253.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = '  -0.05%  '

# Row 24
$ws.Range("D24").Value2 = "'2.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  -0.32%  '

# Row 25
$ws.Range("D25").Value2 = "'2.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = '  -3.73%  '

# Row 26
$ws.Range("D26").Value2 = "'26.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = '  -6.87%  '

# Row 27
$ws.Range("E27").Value2 = '  +0.19%  '

# Row 28
$ws.Range("E28").Value2 = '  +10.39%  '

# Row 29
$ws.Range("E29").Value2 = '  +1.61%  '

# Row 30
$ws.Range("D30").Value2 = "'38.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = '  +4.07%  '

# Row 31
$ws.Range("D31").Value2 = "'6.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = '  +1.61%  '

# Row 32
$ws.Range("D32").Value2 = "'157.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = '  +2.20%  '

# Row 34
$ws.Range("E34").Value2 = '  -2.69%  '

# Row 35
$ws.Range("D35").Value2 = "'0.0794"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = '  -0.56%  '

# Row 36
$ws.Range("D36").Value2 = "'2.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = '  -4.24%  '

# Row 37
$ws.Range("E37").Value2 = '  +2.55%  '

# Row 38
$ws.Range("D38").Value2 = "'18.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = '  -0.78%  '

# Row 39
$ws.Range("E39").Value2 = '  +0.44%  '

# Row 40
$ws.Range("D40").Value2 = "'24.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  +3.80%  '

# Row 41
$ws.Range("E41").Value2 = '  +0.93%  '

# Row 42
$ws.Range("B42").Value2 = 'ApeXProtocol'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").Value2 = "'2.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = '  -2.09%  '

# Row 43
$ws.Range("B43").Value2 = 'RenderToken'
$ws.Range("C43").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value2 = "'3.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = '  +0.65%  '

# Row 44
$ws.Range("D44").Value2 = "'0.0303"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = '  -2.53%  '

# Row 45
$ws.Range("D45").Value2 = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = '  -0.01%  '

# Row 46
$ws.Range("D46").Value2 = '2.046.64'

# Row 47
$ws.Range("D47").Value2 = "'86.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = '  +1.72%  '

# Row 48
$ws.Range("D48").Value2 = "'8.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = '  -3.32%  '

# Row 49
$ws.Range("D49").Value2 = '2.779.56'
$ws.Range("E49").Value2 = '  -1.42%  '

# Row 50
$ws.Range("E50").Value2 = '  +0.51%  '

# Row 51
$ws.Range("D51").Value2 = "'102.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = '  -3.66%  '
